# "threading on left pi"
# Insert a new calibration-data row for the left Pi camera (row 68), update the
# existing data rows with new measurements, and extend the differencing /
# distance formulas that follow the data block down by one row to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Insert a new row at 68 - this shifts the old rows 68:71 (the difference /
# SQRT formulas) down to 69:72, and Excel auto-adjusts all of the relative
# formula references and the shared-formula ranges for us.
$ws.Rows.Item(68).Insert()

# Updated measurement values for the existing data rows (A:C), rows 63-67.
$ws.Range("A63").Value = -3.609
$ws.Range("B63").Value = 1.529
$ws.Range("C63").Value = 35.934

$ws.Range("A64").Value = 0.486
$ws.Range("B64").Value = 1.527
$ws.Range("C64").Value = 40.79

$ws.Range("A65").Value = -3.693
$ws.Range("B65").Value = 1.686
$ws.Range("C65").Value = 47.163

$ws.Range("A66").Value = -10.607
$ws.Range("B66").Value = 1.74
$ws.Range("C66").Value = 48.685

$ws.Range("A67").Value = -10.49
$ws.Range("B67").Value = 1.82
$ws.Range("C67").Value = 62.884

# New row of raw measurement data (row 68, now blank after the insert).
$ws.Range("A68").Value = -9.312
$ws.Range("B68").Value = -4.067
$ws.Range("C68").Value = 62.884

# Extend the difference/distance formula block by one more row (73) so it
# picks up the new data row (68) vs. the previous one (67).
$ws.Range("A73").Formula = "=A68-A67"
$ws.Range("B73").Formula = "=B68-B67"
$ws.Range("C73").Formula = "=C68-C67"
$ws.Range("D73").Formula = "=SQRT((A73^2)+(B73^2)+(C73^2))"

# The stray one-off ratio formulas (24/28 and 48/66) that used to live at the
# end of rows 68 and 71 (now shifted to 69 and 72) are no longer present.
$ws.Range("F69").ClearContents()
$ws.Range("F72").ClearContents()

# Restore the view roughly where it ends up after this edit: scrolled down so
# row 48 is at the top, with the newly added row selected.
$excel.ActiveWindow.ScrollRow = 48
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A72:D73").Select()

Write-Host "done"
